$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row to append: row 34, mirroring the existing weekly price records for
# "Agrícola del Norte S.A. de Arica" - Jengibre (Hortaliza) subset.
$row = 34

$ws.Cells.Item($row, 1).Value = 1
$ws.Cells.Item($row, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item($row, 3).Value = "Arica y Parinacota"
$ws.Cells.Item($row, 4).Value = 45239
$ws.Cells.Item($row, 5).Value = 15
$ws.Cells.Item($row, 6).Value = 100114007
$ws.Cells.Item($row, 7).Value = "Jengibre"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 800
$ws.Cells.Item($row, 11).Value = 16000
$ws.Cells.Item($row, 12).Value = 17000
$ws.Cells.Item($row, 13).Value = 16500
$ws.Cells.Item($row, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item($row, 15).Value = "Perú"
$ws.Cells.Item($row, 16).Value = 1269
$ws.Cells.Item($row, 17).Value = 13
$ws.Cells.Item($row, 18).Value = "Hortaliza"

# Match the date-formatted style already used in column D for other rows.
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item($row - 1, 4).NumberFormat
